# TAP_Q3_2008.xlsx -- "chandra manual annot complete"
# Fill in the manual annotation scores (Clear, Assertive, Cautious,
# Optimistic, Specific, Relevant) for every Q&A row (rows 2-26, cols E:J),
# then restore the reviewer's scroll/zoom/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> [Clear, Assertive, Cautious, Optimistic, Specific, Relevant]
$data = @{
    2  = @(2,2,1,2,1,2)
    3  = @(2,2,1,1,1,2)
    4  = @(2,2,2,2,2,2)
    5  = @(2,1,1,1,1,2)
    6  = @(2,1,2,1,2,2)
    7  = @(2,1,0,1,1,2)
    8  = @(2,1,1,1,1,2)
    9  = @(2,2,1,1,2,2)
    10 = @(2,2,1,1,1,2)
    11 = @(2,2,1,1,2,2)
    12 = @(2,1,1,1,2,2)
    13 = @(2,0,1,2,1,2)
    14 = @(2,0,0,1,0,2)
    15 = @(2,2,1,2,2,2)
    16 = @(2,1,2,1,1,2)
    17 = @(2,1,2,1,1,2)
    18 = @(2,1,2,1,1,2)
    19 = @(2,2,1,2,1,2)
    20 = @(2,0,0,1,0,2)
    21 = @(2,2,2,1,1,2)
    22 = @(2,2,1,1,2,2)
    23 = @(2,1,2,1,1,2)
    24 = @(2,2,1,1,1,2)
    25 = @(2,2,2,2,2,2)
    26 = @(2,1,1,1,1,2)
}

foreach ($r in 2..26) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i  # column E is 5 .. column J is 10
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}

# Restore the reviewer's window state: header row frozen, scrolled down to
# row 21, zoomed to 85%, with E27 selected.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.Goto($ws.Range("A21"), $true)
$excel.ActiveWindow.Zoom = 85
$ws.Range("E27").Select()
